# remove gas loss factor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the row first (mirrors how the edit was made in Excel - whole row
# selected then deleted), then delete it so everything below shifts up.
$ws.Rows.Item(5).Select()
$ws.Rows.Item(5).Delete()

# Leave the selection where Excel would land after a row delete.
$ws.Range("A5:XFD5").Select()

# Update the defined names: GasLossFactor pointed straight at the deleted
# cell (B5), so it becomes a broken reference; everything else that lived
# below row 5 shifts up by one row.
$wb.Names.Item("GasLossFactor").RefersTo = "=Sheet1!#REF!"
$wb.Names.Item("ExplorationFee").RefersTo = "=Sheet1!`$B`$12"
$wb.Names.Item("GasProd").RefersTo = "=Sheet1!`$B`$7:`$J`$7"
$wb.Names.Item("OilProd").RefersTo = "=Sheet1!`$B`$6:`$J`$6"
$wb.Names.Item("PostExplorationFee").RefersTo = "=Sheet1!`$B`$13"
